# Natmi following Dr Hou advice
# Rebuild the LR-pair result rows (2-7) on Sheet1 with the recomputed
# specificity figures, and extend the table with the new FAPs/sCs x ECs
# combinations that came out of the updated pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D are the cluster/gene labels, columns E-T are the numeric
# NATMI edge-scoring metrics (detection rates, expression values and
# specificities).
$rows = @(
    @{ Row = 2;  A = "FAPs"; B = "Btc"; C = "Erbb2"; D = "ECs";  V = @(2, 0.6666666666666666, 0.191714, 0.575142, 0.09369188973541917, 0.09369188973541917, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 0.3297809269986667, 2.968028342988, 0.01837725526640716, 0.01837725526640716) },
    @{ Row = 3;  A = "FAPs"; B = "Btc"; C = "Erbb2"; D = "FAPs"; V = @(2, 0.6666666666666666, 0.191714, 0.575142, 0.09369188973541917, 0.09369188973541917, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 0.7815007768413333, 7.033506991572001, 0.04354963580706639, 0.04354963580706639) },
    @{ Row = 4;  A = "FAPs"; B = "Btc"; C = "Erbb2"; D = "sCs";  V = @(2, 0.6666666666666666, 0.191714, 0.575142, 0.09369188973541917, 0.09369188973541917, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 0.570024769912, 5.130222929208, 0.03176499866194562, 0.03176499866194562) },
    @{ Row = 5;  A = "sCs";  B = "Btc"; C = "Erbb2"; D = "ECs";  V = @(3, 1, 1.854503666666667, 5.563511, 0.9063081102645809, 0.9063081102645809, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 3.190064044961556, 28.710576404654, 0.1777683803729586, 0.1777683803729586) },
    @{ Row = 6;  A = "sCs";  B = "Btc"; C = "Erbb2"; D = "FAPs"; V = @(3, 1, 1.854503666666667, 5.563511, 0.9063081102645809, 0.9063081102645809, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 7.559677729091778, 68.037099561826, 0.4212679266313497, 0.4212679266313497) },
    @{ Row = 7;  A = "sCs";  B = "Btc"; C = "Erbb2"; D = "sCs";  V = @(3, 1, 1.854503666666667, 5.563511, 0.9063081102645809, 0.9063081102645809, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 5.514010588129334, 49.626095293164, 0.3072718032602726, 0.3072718032602726) }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D

    for ($i = 0; $i -lt $r.V.Count; $i++) {
        $colIndex = 5 + $i
        $ws.Cells.Item($rowIndex, $colIndex).Value = $r.V[$i]
    }
}
